# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Fri Sep 27 19:51:34 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.755.61'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.699.35'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('E9').Value = '  +5.02%  '
$ws.Range('E10').Value = '  +4.61%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '30.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('E14').Value = '  +9.30%  '
$ws.Range('D15').Value = '3.185.17'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = '65.626.04'
$ws.Range('D17').Value = '2.709.29'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('E25').Value = '  +12.26%  '
$ws.Range('E26').Value = '  -4.78%  '
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('E28').Value = '  +3.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('E30').Value = '  +3.16%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '537.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.54%  '
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.71'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '167.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0615'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.24%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.660'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0266'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.86%  '
$ws.Range('E51').Value = '  -0.30%  '
